$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '3.66%'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '27.63'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-8.29%'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.220'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '1.25%'

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05873'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '1.91%'

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '0.73%'

$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.219'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-1.76%'

$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8634'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '1.61%'

$ws.Range("B9").Value = 'FTXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9621'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '11.92%'

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '2.03%'

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07162'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '1.10%'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03180'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-1.69%'

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09209'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-1.63%'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001555'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '1.63%'

$ws.Range("B15").Value = 'One'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0006041'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-94.11%'

$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005795'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-1.76%'

$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.501'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-1.27%'

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '0.22%'

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '1.18%'

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.03478'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '1.63%'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1307'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-0.69%'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.551'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '2.09%'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04161'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '1.26%'

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-2.09%'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001225'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '0.25%'

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004795'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '15.31%'

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '0.04%'

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '1.22%'

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '1.53%'

$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1101'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '2.86%'

$ws.Range("B42").Value = 'KickToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.003820'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-32.95%'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002344'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-4.69%'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01069'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '4.78%'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005239'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-4.41%'

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '0.04%'

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1000'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '40.90%'

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002101'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '0.04%'

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0002001'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.04%'
